# TRE-95 (OR-550-BE): Fix branch revenue
# Insert a new column before the existing "Tổng xu" column for the new
# metric "Tổng giảm giá đơn dưới 2.000đ", shifting the remaining totals
# columns (Tổng xu, Tổng doanh thu, Tổng doanh thu tiền mặt, Tổng doanh
# thu chuyển khoản ngân hàng, Tổng doanh thu ví nội bộ) one column right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at G - this shifts the old G:K headers to H:L and
# extends the merged title cell A1:K1 to A1:L1 automatically.
$ws.Range("G1").EntireColumn.Insert()

# Give the new column the same width as its neighbours.
$ws.Range("G1").EntireColumn.ColumnWidth = $ws.Range("F1").EntireColumn.ColumnWidth

# Copy the header formatting (fill/border/font/alignment) from an
# existing header cell onto the new one, then set its text.
$ws.Range("H8").Copy()
$ws.Range("G8").PasteSpecial(-4122)
$ws.Range("G8").Value = "Tổng giảm giá đơn dưới 2.000đ"
$excel.CutCopyMode = $false

# Restore the active selection.
[void]$ws.Range("G16").Select()
